# Rename the "12-3" sheet to "12-03" (also fixes the _FilterDatabase
# defined name that refers to the sheet by its quoted name), then move
# the active selection on that sheet to F17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("12-3")
$ws.Name = "12-03"

$ws.Activate()
$ws.Range("F17").Select()
